# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on the
# per-language detail sheets to reflect the newly generated report run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-12 00:46:40"
$zhcn.Range("H2").Value = "2016-03-12 00:46:57"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-12 00:46:43"
$dede.Range("H2").Value = "2016-03-12 00:47:02"
